$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "(188.0, 211.0, 220.0)"
$ws.Range("C2").Value = "(182.0, 216.0, 45.0)"
$ws.Range("D2").Value = "(245.0, 150.0, 5.0)"
$ws.Range("E2").Value = "(240.0, 192.0, 223.0)"
$ws.Range("F2").Value = "(186.0, 101.0, 78.0)"
$ws.Range("G2").Value = "(90.5, 97.5, 87.0)"
$ws.Range("H2").Value = "(230.0, 229.0, 225.0)"
$ws.Range("I2").Value = "(171.0, 158.0, 31.0)"
$ws.Range("B3").Value = "(237.0, 253.0, 254.0)"
$ws.Range("C3").Value = "(128.0, 199.5, 122.0)"
$ws.Range("E3").Value = "(255.0, 157.0, 255.0)"
$ws.Range("F3").Value = "(254.0, 80.0, 60.0)"
$ws.Range("G3").Value = "(113.0, 110.0, 108.0)"
$ws.Range("H3").Value = "(228.0, 228.0, 227.0)"
$ws.Range("I3").Value = "(254.0, 255.0, 53.0)"
$ws.Range("B4").Value = "(227.0, 235.0, 242.0)"
$ws.Range("G4").Value = "(136.0, 130.0, 123.0)"
$ws.Range("H4").Value = "(235.0, 234.0, 230.0)"
$ws.Range("I4").Value = "(126.0, 86.0, 32.0)"
$ws.Range("B5").Value = "(226.0, 235.0, 241.0)"
$ws.Range("G5").Value = "(137.0, 132.0, 124.0)"
$ws.Range("H5").Value = "(233.0, 232.0, 234.0)"
$ws.Range("B6").Value = "(169.0, 217.0, 243.0)"
$ws.Range("G6").Value = "(138.0, 132.0, 118.0)"
$ws.Range("H6").Value = "(237.0, 236.0, 238.0)"
$ws.Range("B7").Value = "(167.0, 208.0, 244.0)"
$ws.Range("G7").Value = "(136.0, 127.5, 115.0)"
$ws.Range("H7").Value = "(225.0, 225.0, 233.0)"
$ws.Range("B8").Value = "(224.0, 235.0, 240.0)"
$ws.Range("B9").Value = "(184.0, 217.0, 250.0)"
$ws.Range("B10").Value = "(35.0, 158.0, 235.0)"
$ws.Range("B11").Value = "(41.0, 144.0, 211.0)"
$ws.Range("B12").Value = "(242.0, 254.0, 254.0)"
$ws.Range("B13").Value = "(226.0, 245.0, 252.0)"
$ws.Range("B14").Value = "(241.0, 253.0, 253.0)"
$ws.Range("B15").Value = "(239.0, 254.0, 253.0)"
$ws.Range("B16").Value = "(202.0, 226.0, 238.0)"
$ws.Range("B17").Value = "(197.0, 227.0, 232.0)"
$ws.Range("B18").Value = "(229.0, 239.0, 248.0)"
$ws.Range("B19").Value = "(247.0, 255.0, 254.0)"
